# Add a "Serial No" column at the front of the Teams Data sheet, and fix a
# few team-roster / powerup-order typos in the Users / Powerups columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column A to hold the serial number.
$ws.Columns.Item(1).Insert()

# Header for the new column.
$ws.Range("A1").Value = "Serial No"

# Column widths: new A (Serial No) = 10, new B (Team Name) = 20,
# new C (Users) = 50, new D:H (Phase/Powerups cols) = 20, new I (Credit Card No) = 30.
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 50
$ws.Range("D1:H1").EntireColumn.ColumnWidth = 20
$ws.Columns.Item(9).ColumnWidth = 30

# Fill the serial numbers 1..28 for the 28 data rows (rows 2-29).
for ($i = 2; $i -le 29; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Fix roster typos now that Users moved from column B to column C.
$ws.Range("C2").Value = "Ankush Gautam, Vaibhav Srivastva, Mahi, Aishlee Joshi"
$ws.Range("C22").Value = "Aayushman, Madhav Gaba, Saksham Katna, hemant"

# Fix the Powerups order value (now column H) for the IMPOSTORS row.
$ws.Range("H24").Value = "6, 7, 8"

Write-Output "Serial No column inserted and roster/powerup fixes applied."
